# Insert a new data row at row 301 (pushing the existing rows 301..418 down
# to 302..419) and populate the new row with the reported price observation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("301:301").Insert()

$ws.Cells.Item(301, 1).Value = 10
$ws.Cells.Item(301, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(301, 3).Value = "La Araucanía"
$ws.Cells.Item(301, 4).Value = 44510
$ws.Cells.Item(301, 5).Value = 9
$ws.Cells.Item(301, 6).Value = 100112003
$ws.Cells.Item(301, 7).Value = "Ajo"
$ws.Cells.Item(301, 8).Value = "Chino"
$ws.Cells.Item(301, 9).Value = "Primera"
$ws.Cells.Item(301, 10).Value = 200
$ws.Cells.Item(301, 11).Value = 20000
$ws.Cells.Item(301, 12).Value = 21000
$ws.Cells.Item(301, 13).Value = 20500
$ws.Cells.Item(301, 14).Value = "`$/caja 10 kilos"
$ws.Cells.Item(301, 15).Value = "China"
$ws.Cells.Item(301, 16).Value = 2050
$ws.Cells.Item(301, 17).Value = 10
$ws.Cells.Item(301, 18).Value = "Hortaliza"
